# Update the "handback-status" report timestamps to reflect the latest
# generation run (commit: "Generate Report for Handback").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G)
$wsOverview.Range("G2").Value = "2016-08-26 15:20:02"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$wsZhCn.Range("H2").Value = "2016-08-26 15:19:56"
$wsZhCn.Range("K2").Value = "2016-08-26 15:20:35"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
$wsDeDe.Range("H2").Value = "2016-08-26 15:20:02"
$wsDeDe.Range("K2").Value = "2016-08-26 15:20:43"
